$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 237
$ws.Range("I2").Value = 626
$ws.Range("J2").Value = 2676
$ws.Range("K2").Value = 12
$ws.Range("L2").Value = 742
$ws.Range("M2").Value = 46
$ws.Range("N2").Value = 444
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 7
$ws.Range("Q2").Value = 6
$ws.Range("R2").Value = 26
$ws.Range("S2").Value = 301
$ws.Range("T2").Value = 517
$ws.Range("U2").Value = 36
$ws.Range("V2").Value = 4214
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 4219
$ws.Range("Y2").Value = 7
$ws.Range("Z2").Value = 78
$ws.Range("AA2").Value = 25
